$d = $word.ActiveDocument

# The document has a "Heading 2" paragraph containing just "ACT" (the book
# code) immediately followed by a duplicate, italicized "Kisah Para Rasul"
# paragraph (the book's full title, shown again a few paragraphs later as the
# real "Heading 2" section title). That duplicate paragraph is redundant and
# needs to be removed entirely -- including its paragraph mark -- so the
# "ACT" heading paragraph is followed directly by what used to come after the
# duplicate.

$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count -and -not $found; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)

    if ($text -eq "ACT") {
        $dup = $p.Next()
        $dupText = $dup.Range.Text.TrimEnd([char]13, [char]7)

        if ($dupText -eq "Kisah Para Rasul" -and $dup.Range.Italic) {
            $dup.Range.Delete()
            $found = $true
        }
    }
}
